$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stores Price (D) / Volume 1h (E) as literal text (e.g. "303.37",
# "1.29%"), not as numbers. Prefixing the literal with a single quote forces
# Excel to keep the entry as text instead of auto-converting it to a
# number/percentage, which mirrors how the source data was authored.
#
# Row => @{ D = newPrice; E = newVolume }
$updates = @{
    2  = @{ D = "303.37";        E = "1.29%" }
    3  = @{ D = "32.91";         E = "5.08%" }
    4  = @{ D = "4.948";         E = "-3.10%" }
    5  = @{ D = "0.07841";       E = "-1.55%" }
    6  = @{ D = "1.997";         E = "-15.08%" }
    7  = @{ D = "7.845";         E = "0.98%" }
    8  = @{ D = "3.805";         E = "-1.49%" }
    9  = @{ D = "0.9264";        E = "0.39%" }
    10 = @{ D = "0.1756";        E = "1.16%" }
    11 = @{ D = "0.07813";       E = "3.04%" }
    12 = @{ D = "0.08676";       E = "-7.13%" }
    13 = @{ D = "0.03142";       E = "3.20%" }
    14 = @{ E = "0.16%" }
    15 = @{ E = "0.45%" }
    16 = @{ D = "0.005921";      E = "-1.67%" }
    17 = @{ D = "3.466";         E = "-0.55%" }
    18 = @{ D = "2.155";         E = "-5.00%" }
    19 = @{ E = "1.15%" }
    20 = @{ E = "-1.30%" }
    21 = @{ D = "4.312";         E = "9.85%" }
    22 = @{ E = "17.12%" }
    23 = @{ E = "-1.46%" }
    24 = @{ E = "-1.98%" }
    25 = @{ D = "0.004450";      E = "-0.69%" }
    26 = @{ D = "0.0001250";     E = "4.18%" }
    39 = @{ E = "-1.03%" }
    40 = @{ D = "0.04794";       E = "3.51%" }
    41 = @{ D = "0.007501";      E = "7.54%" }
    42 = @{ E = "-0.03%" }
    43 = @{ D = "0.002339";      E = "6.86%" }
    44 = @{ D = "0.01172";       E = "14.04%" }
    45 = @{ D = "0.00006242";    E = "-0.54%" }
    46 = @{ D = "0.00000000750"; E = "0.02%" }
    47 = @{ E = "-61.12%" }
    48 = @{ D = "0.8206";        E = "9.89%" }
    49 = @{ D = "0.00002100";    E = "0.02%" }
    50 = @{ D = "0.0002000";     E = "0.02%" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $ws.Range("D$row").Value = "'" + $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Range("E$row").Value = "'" + $rowData["E"]
    }
}
